# Add row 27 (new company news item: Farrel Pomini) to the bottom of the
# news table on Sheet1, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 27

# Helper: write a plain text value into a cell without Excel's
# autodetection turning date-like strings (e.g. "2025-11-21") into date
# serial numbers, and without leaving the cell's style pointed at a
# non-default "text" format.
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Helper: write an "empty string" cell (present in the sheet, but with no
# visible text) matching the inline-string-with-empty-content cells used
# throughout the rest of the sheet (e.g. H4, H7, L26, ...).
function Set-EmptyTextCell($cell) {
    $cell.NumberFormat = "@"
    $cell.Value = "Farrel-Pomini-EMPTY-PLACEHOLDER"
    $cell.Value = ""
    $cell.Style = "Normal"
}

Set-TextCell ($ws.Cells.Item($newRow, 1)) "26"
$ws.Cells.Item($newRow, 1).Value = 26

Set-TextCell ($ws.Cells.Item($newRow, 2)) "Farrel Pomini"
Set-EmptyTextCell ($ws.Cells.Item($newRow, 3))
Set-TextCell ($ws.Cells.Item($newRow, 4)) "2025-11-21"
Set-TextCell ($ws.Cells.Item($newRow, 5)) "Farrel Pomini Pioneers a Greener Tomorrow Through Innovation"
Set-TextCell ($ws.Cells.Item($newRow, 6)) "Farrel Pomini announced new advancements and sustainability initiatives in biopolymer processing, reinforcing its leadership in continuous mixing technology."
Set-TextCell ($ws.Cells.Item($newRow, 7)) "Product Launch"
Set-EmptyTextCell ($ws.Cells.Item($newRow, 8))
Set-TextCell ($ws.Cells.Item($newRow, 9)) "https://markets.financialcontent.com/wral/article/tokenring-2025-11-21-farrel-pomini-pioneers-a-greener-tomorrow-through-relentless-innovation-in-manufacturing"
Set-TextCell ($ws.Cells.Item($newRow, 10)) "2025-W47"
Set-TextCell ($ws.Cells.Item($newRow, 11)) "Perplexity Rev2"
Set-EmptyTextCell ($ws.Cells.Item($newRow, 12))
Set-TextCell ($ws.Cells.Item($newRow, 13)) "Yes"
Set-TextCell ($ws.Cells.Item($newRow, 14)) "2025-11-21-farrel-pomini-pioneers-a-greener-tomorrow-through-innovation.md"
